# Update cryptocurrency price (column D) and 1h volume-change (column E) values
# to reflect refreshed data from the GitHub Actions scheduled scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.863.80'   # D2: '26.873.86' -> '26.863.80'
$ws.Cells.Item(2, 5).Value = '  +0.14%  '   # E2: '  +0.21%  ' -> '  +0.14%  '

$ws.Cells.Item(3, 4).Value = '1.639.94'   # D3: '1.639.83' -> '1.639.94'
$ws.Cells.Item(3, 5).Value = '  -0.18%  '   # E3: '  -0.24%  ' -> '  -0.18%  '

$ws.Cells.Item(4, 5).Value = '  -0.57%  '   # E4: '  -0.52%  ' -> '  -0.57%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '216.77'   # D5: '216.82' -> '216.77'
$ws.Cells.Item(5, 5).Value = '  -0.84%  '   # E5: '  -0.77%  ' -> '  -0.84%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.512'   # D6: '0.511' -> '0.512'
$ws.Cells.Item(6, 5).Value = '  +2.04%  '   # E6: '  +1.84%  ' -> '  +2.04%  '

$ws.Cells.Item(7, 5).Value = '  -0.53%  '   # E7: '  -0.49%  ' -> '  -0.53%  '

$ws.Cells.Item(8, 5).Value = '  +1.69%  '   # E8: '  +1.75%  ' -> '  +1.69%  '

$ws.Cells.Item(9, 5).Value = '  +0.35%  '   # E9: '  +0.43%  ' -> '  +0.35%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '19.87'   # D10: '19.88' -> '19.87'
$ws.Cells.Item(10, 5).Value = '  +3.13%  '   # E10: '  +3.35%  ' -> '  +3.13%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0847'   # D11: '0.0848' -> '0.0847'
$ws.Cells.Item(11, 5).Value = '  -0.25%  '   # E11: '  -0.07%  ' -> '  -0.25%  '

$ws.Cells.Item(12, 4).Value = '1.869.37'   # D12: '1.869.17' -> '1.869.37'
$ws.Cells.Item(12, 5).Value = '  -0.20%  '   # E12: '  -0.18%  ' -> '  -0.20%  '

$ws.Cells.Item(13, 4).Value = '1.636.02'   # D13: '1.638.46' -> '1.636.02'
$ws.Cells.Item(13, 5).Value = '  -0.42%  '   # E13: '  -0.01%  ' -> '  -0.42%  '

$ws.Cells.Item(14, 5).Value = '  -0.73%  '   # E14: '  -0.71%  ' -> '  -0.73%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.530'   # D15: '0.531' -> '0.530'
$ws.Cells.Item(15, 5).Value = '  +0.80%  '   # E15: '  +0.92%  ' -> '  +0.80%  '

$ws.Cells.Item(16, 5).Value = '  +2.99%  '   # E16: '  +3.11%  ' -> '  +2.99%  '

$ws.Cells.Item(17, 4).Value = '26.861.03'   # D17: '26.868.04' -> '26.861.03'
$ws.Cells.Item(17, 5).Value = '  +0.06%  '   # E17: '  +0.18%  ' -> '  +0.06%  '

$ws.Cells.Item(18, 5).Value = '  -0.63%  '   # E18: '  -0.64%  ' -> '  -0.63%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '219.32'   # D19: '219.90' -> '219.32'
$ws.Cells.Item(19, 5).Value = '  +1.51%  '   # E19: '  +2.07%  ' -> '  +1.51%  '

$ws.Cells.Item(20, 5).Value = '  -0.54%  '   # E20: '  -0.53%  ' -> '  -0.54%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '6.84'   # D21: '6.86' -> '6.84'
$ws.Cells.Item(21, 5).Value = '  +3.20%  '   # E21: '  +3.69%  ' -> '  +3.20%  '

$ws.Cells.Item(22, 5).Value = '  +0.45%  '   # E22: '  +0.41%  ' -> '  +0.45%  '

$ws.Cells.Item(23, 5).Value = '  +3.41%  '   # E23: '  +3.68%  ' -> '  +3.41%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '9.15'   # D24: '9.16' -> '9.15'
$ws.Cells.Item(24, 5).Value = '  -0.66%  '   # E24: '  -0.47%  ' -> '  -0.66%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '147.18'   # D25: '147.15' -> '147.18'
$ws.Cells.Item(25, 5).Value = '  -0.33%  '   # E25: '  -0.40%  ' -> '  -0.33%  '

$ws.Cells.Item(26, 5).Value = '  -0.47%  '   # E26: '  -0.54%  ' -> '  -0.47%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '7.34'   # D27: '7.35' -> '7.34'
$ws.Cells.Item(27, 5).Value = '  +2.75%  '   # E27: '  +2.98%  ' -> '  +2.75%  '

$ws.Cells.Item(28, 5).Value = '  +0.53%  '   # E28: '  +0.28%  ' -> '  +0.53%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '15.79'   # D29: '15.80' -> '15.79'
$ws.Cells.Item(29, 5).Value = '  +0.24%  '   # E29: '  +0.48%  ' -> '  +0.24%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.0503'   # D30: '0.0504' -> '0.0503'
$ws.Cells.Item(30, 5).Value = '  -1.41%  '   # E30: '  -1.21%  ' -> '  -1.41%  '

$ws.Cells.Item(31, 5).Value = '  -1.06%  '   # E31: '  -0.87%  ' -> '  -1.06%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.33'   # D32: '3.34' -> '3.33'
$ws.Cells.Item(32, 5).Value = '  -1.54%  '   # E32: '  -1.32%  ' -> '  -1.54%  '

$ws.Cells.Item(33, 5).Value = '  +0.46%  '   # E33: '  +0.50%  ' -> '  +0.46%  '

$ws.Cells.Item(34, 5).Value = '  +1.39%  '   # E34: '  +1.32%  ' -> '  +1.39%  '

$ws.Cells.Item(35, 4).Value = '1.267.45'   # D35: '1.265.50' -> '1.267.45'
$ws.Cells.Item(35, 5).Value = '  -0.12%  '   # E35: '  -0.26%  ' -> '  -0.12%  '

$ws.Cells.Item(36, 5).Value = '  -0.24%  '   # E36: '  -0.14%  ' -> '  -0.24%  '

$ws.Cells.Item(37, 5).Value = '  +1.61%  '   # E37: '  +2.06%  ' -> '  +1.61%  '

$ws.Cells.Item(38, 5).Value = '  +0.11%  '   # E38: '  +0.34%  ' -> '  +0.11%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.833'   # D39: '0.834' -> '0.833'
$ws.Cells.Item(39, 5).Value = '  +1.94%  '   # E39: '  +2.08%  ' -> '  +1.94%  '

$ws.Cells.Item(40, 5).Value = '  -0.48%  '   # E40: '  -0.45%  ' -> '  -0.48%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.810'   # D41: '0.811' -> '0.810'
$ws.Cells.Item(41, 5).Value = '  +0.71%  '   # E41: '  +0.98%  ' -> '  +0.71%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '5.40'   # D42: '5.39' -> '5.40'
$ws.Cells.Item(42, 5).Value = '  +0.88%  '   # E42: '  +0.95%  ' -> '  +0.88%  '

$ws.Cells.Item(43, 4).Value = '1.780.55'   # D43: '1.779.39' -> '1.780.55'
$ws.Cells.Item(43, 5).Value = '  -0.16%  '   # E43: '  -0.13%  ' -> '  -0.16%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '61.89'   # D44: '61.87' -> '61.89'
$ws.Cells.Item(44, 5).Value = '  +1.04%  '   # E44: '  +0.70%  ' -> '  +1.04%  '

$ws.Cells.Item(45, 5).Value = '  -1.56%  '   # E45: '  -1.48%  ' -> '  -1.56%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '91.72'   # D46: '91.82' -> '91.72'
$ws.Cells.Item(46, 5).Value = '  -1.13%  '   # E46: '  -1.10%  ' -> '  -1.13%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.58'   # D47: '1.59' -> '1.58'
$ws.Cells.Item(47, 5).Value = '  -1.38%  '   # E47: '  -0.97%  ' -> '  -1.38%  '

$ws.Cells.Item(48, 5).Value = '  +2.93%  '   # E48: '  +1.12%  ' -> '  +2.93%  '

$ws.Cells.Item(49, 5).Value = '  -0.58%  '   # E49: '  -0.42%  ' -> '  -0.58%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '7.63'   # D50: '7.65' -> '7.63'
$ws.Cells.Item(50, 5).Value = '  +0.78%  '   # E50: '  +1.52%  ' -> '  +0.78%  '

$ws.Cells.Item(51, 5).Value = '  -0.54%  '   # E51: '  -0.31%  ' -> '  -0.54%  '
